# "Continuation intégration structure stcReleve"
#
# The item #27 (row 32, "SAUVEGARDE" / "Récupérer le répertoire par ftp...")
# moves from status "En cours" to status "Clos". Because the sheet's
# AutoFilter only shows rows whose Statut is "Ouvert" or "En cours", this
# row (and row 24, item #19, whose Statut was already "Clos") end up
# filtered out / hidden. The dependent COUNTIF() summary cells (E2/F2)
# recalculate automatically. Finally the frozen-pane scroll position is
# nudged down a few rows (B28 -> B31) to keep the newly-hidden rows out of
# view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status of item #27 (row 32) from "En cours" to "Clos".
$ws.Range("G32").Value = "Clos"

# Re-applying the filter (Statut in {Ouvert, En cours}) hides rows that no
# longer qualify: row 24 (already Clos) and row 32 (now Clos too).
$ws.Rows.Item(24).Hidden = $true
$ws.Rows.Item(32).Hidden = $true

# Scroll the frozen window down a bit so the view starts around row 31
# instead of row 28.
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 2
